$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 hold two separate sightings that were recorded in swapped
# order. Columns A, B, D, E, F, G, H, M, Q, R differ between the rows and
# need to be swapped; every other column already matches between the two
# rows, so leaving them untouched is a no-op.
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range("$col" + "2")
    $cell3 = $ws.Range("$col" + "3")

    # Value2 round-trips numbers/strings cleanly through this COM host;
    # plain Value can box into a Variant wrapper that doesn't re-set well.
    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
